$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 19
$ws.Range("G19").Value = 1.62
$ws.Range("I19").Value = 6.25
$ws.Range("Q19").Value = 2.2
$ws.Range("R19").Value = 1.65
$ws.Range("U19").Value = 2.2
$ws.Range("V19").Value = 1.62
$ws.Range("X19").Value = 6.5
$ws.Range("AH19").Value = 13
$ws.Range("AI19").Value = 29
$ws.Range("AK19").Value = 67
$ws.Range("AL19").Value = 51
$ws.Range("AN19").Value = 3.4
$ws.Range("AO19").Value = 8.5
$ws.Range("AY19").Value = 34

# Row 48
$ws.Range("G48").Value = 3.6
$ws.Range("I48").Value = 2.05
$ws.Range("L48").Value = 2.63
$ws.Range("U48").Value = 1.57
$ws.Range("V48").Value = 2.25
$ws.Range("X48").Value = 19
$ws.Range("Y48").Value = 12
$ws.Range("AG48").Value = 126
$ws.Range("AH48").Value = 9.5
$ws.Range("AO48").Value = 17

# Row 50
$ws.Range("G50").Value = 3.1
$ws.Range("H50").Value = 2.95
$ws.Range("I50").Value = 2.32
$ws.Range("J50").Value = 3.65
$ws.Range("L50").Value = 2.87
$ws.Range("M50").Value = 1.08
$ws.Range("N50").Value = 7.5
$ws.Range("O50").Value = 1.36
$ws.Range("P50").Value = 2.7
$ws.Range("T50").Value = 2.52
$ws.Range("U50").Value = 1.75
$ws.Range("V50").Value = 1.85
$ws.Range("W50").Value = 8.5
$ws.Range("Y50").Value = 11
$ws.Range("AA50").Value = 30
$ws.Range("AB50").Value = 37
$ws.Range("AD50").Value = 5.8
$ws.Range("AE50").Value = 14
$ws.Range("AF50").Value = 70
$ws.Range("AG50").Value = 600
$ws.Range("AI50").Value = 11.25
$ws.Range("AJ50").Value = 9
$ws.Range("AK50").Value = 25
$ws.Range("AL50").Value = 20
$ws.Range("AM50").Value = 30
$ws.Range("AN50").Value = 5
$ws.Range("AO50").Value = 17
$ws.Range("AP50").Value = 23
$ws.Range("AQ50").Value = 80
$ws.Range("AR50").Value = 120
$ws.Range("AS50").Value = 300
$ws.Range("AT50").Value = 2.47
$ws.Range("AU50").Value = 6.6
$ws.Range("AV50").Value = 55
$ws.Range("AX50").Value = 4.2
$ws.Range("AY50").Value = 12
$ws.Range("AZ50").Value = 19
$ws.Range("BA50").Value = 50
$ws.Range("BB50").Value = 80
